$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Worksheet references
# ------------------------------------------------------------------
$wsTeam      = $wb.Worksheets.Item("ProjectTeam")
$wsProduct   = $wb.Worksheets.Item("Product Backlog")
$wsSprint    = $wb.Worksheets.Item("Sprint Backlog")
$wsBurndown  = $wb.Worksheets.Item("BurndownChart")

# ------------------------------------------------------------------
# Zoom level changed from 110% to 160% on every sheet
# (zoom is a per-window/active-sheet setting, so each sheet has to be
#  activated before the window zoom is changed)
# ------------------------------------------------------------------
$wsTeam.Activate()
$excel.ActiveWindow.Zoom = 160

$wsProduct.Activate()
$excel.ActiveWindow.Zoom = 160

$wsSprint.Activate()
$excel.ActiveWindow.Zoom = 160

$wsBurndown.Activate()
$excel.ActiveWindow.Zoom = 160

# ------------------------------------------------------------------
# Column widths were slightly reduced (best effort - the COM width
# model only supports 1/6 character increments)
# ------------------------------------------------------------------
$wsTeam.Columns.Item(1).ColumnWidth = 7.666666666666667

$wsProduct.Columns.Item(1).ColumnWidth = 7.666666666666667
$wsProduct.Columns.Item(2).ColumnWidth = 21.333333333333332
$wsProduct.Columns.Item(3).ColumnWidth = 33.5

$wsSprint.Columns.Item(1).ColumnWidth = 7.666666666666667
$wsSprint.Columns.Item(2).ColumnWidth = 7.666666666666667
$wsSprint.Columns.Item(3).ColumnWidth = 40.0
$wsSprint.Columns.Item(4).ColumnWidth = 7.666666666666667
$wsSprint.Columns.Item(5).ColumnWidth = 11.166666666666666
$wsSprint.Columns.Item(6).ColumnWidth = 7.666666666666667
$wsSprint.Columns.Item(7).ColumnWidth = 13.666666666666666

$wsBurndown.Columns.Item(1).ColumnWidth = 7.666666666666667

# ------------------------------------------------------------------
# Sprint Backlog data updates
# ------------------------------------------------------------------

# Row 11: reviewer/owner reassigned from Hannes/Nic to Ken/Joel
$wsSprint.Range("F11").Value = "Ken"
$wsSprint.Range("G11").Value = "Joel"

# Row 12: remaining-effort value added
$wsSprint.Range("J12").Value = 4

# Row 13: remaining-effort value added
$wsSprint.Range("J13").Value = 7.5

# Row 14: owner/reviewer swapped back to Hannes/Nic, remaining-effort added
$wsSprint.Range("F14").Value = "Hannes"
$wsSprint.Range("G14").Value = "Nic"
$wsSprint.Range("J14").Value = 4

# Row 15: remaining-effort value added
$wsSprint.Range("J15").Value = 2

# Row 16: remaining-effort value added
$wsSprint.Range("J16").Value = 4

# New Sprint 3 rows
$wsSprint.Range("A17").Value = 3.1
$wsSprint.Range("B17").Value = 3
$wsSprint.Range("B18").Value = 3
$wsSprint.Range("B19").Value = 3
$wsSprint.Range("B20").Value = 3
$wsSprint.Range("B21").Value = 3
$wsSprint.Range("B22").Value = 3

# Selection on Sprint Backlog moved from J12 to F12
[void]$wsSprint.Range("F12").Select()

# Restore Sprint Backlog as the active / selected tab
$wsSprint.Activate()
